$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$nl = [char]11

$t.Cell(1, 1).Range.Text = "93 x 69" + $nl + "  6    9" + $nl + "  ----" + $nl + "9|    |" + $nl + "3|    |"
$t.Cell(1, 2).Range.Text = "92 x 70" + $nl + "  7    0" + $nl + "  ----" + $nl + "9|    |" + $nl + "2|    |"
$t.Cell(1, 3).Range.Text = "28 x 37" + $nl + "  3    7" + $nl + "  ----" + $nl + "2|    |" + $nl + "8|    |"
$t.Cell(2, 1).Range.Text = "40 x 15" + $nl + "  1    5" + $nl + "  ----" + $nl + "4|    |" + $nl + "0|    |"
$t.Cell(2, 2).Range.Text = "39 x 58" + $nl + "  5    8" + $nl + "  ----" + $nl + "3|    |" + $nl + "9|    |"
$t.Cell(2, 3).Range.Text = "69 x 43" + $nl + "  4    3" + $nl + "  ----" + $nl + "6|    |" + $nl + "9|    |"
$t.Cell(3, 1).Range.Text = "92 x 19" + $nl + "  1    9" + $nl + "  ----" + $nl + "9|    |" + $nl + "2|    |"
$t.Cell(3, 2).Range.Text = "73 x 38" + $nl + "  3    8" + $nl + "  ----" + $nl + "7|    |" + $nl + "3|    |"
$t.Cell(3, 3).Range.Text = "47 x 41" + $nl + "  4    1" + $nl + "  ----" + $nl + "4|    |" + $nl + "7|    |"
$t.Cell(4, 1).Range.Text = "34 x 81" + $nl + "  8    1" + $nl + "  ----" + $nl + "3|    |" + $nl + "4|    |"
$t.Cell(4, 2).Range.Text = "53 x 81" + $nl + "  8    1" + $nl + "  ----" + $nl + "5|    |" + $nl + "3|    |"
$t.Cell(4, 3).Range.Text = "55 x 77" + $nl + "  7    7" + $nl + "  ----" + $nl + "5|    |" + $nl + "5|    |"
$t.Cell(5, 1).Range.Text = "59 x 66" + $nl + "  6    6" + $nl + "  ----" + $nl + "5|    |" + $nl + "9|    |"
$t.Cell(5, 2).Range.Text = "38 x 75" + $nl + "  7    5" + $nl + "  ----" + $nl + "3|    |" + $nl + "8|    |"
$t.Cell(5, 3).Range.Text = "88 x 14" + $nl + "  1    4" + $nl + "  ----" + $nl + "8|    |" + $nl + "8|    |"
